$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.07026449749202
$ws.Range("D2").Value = 1.071470147854364
$ws.Range("E2").Value = 1.070553816562493
$ws.Range("F2").Value = 1.082171881303255
$ws.Range("I2").Value = 1.049378931519855
$ws.Range("J2").Value = 1.075195032521656
$ws.Range("K2").Value = 1.074167235954113
$ws.Range("L2").Value = 1.073253345669517
$ws.Range("M2").Value = 1.084840791683628
$ws.Range("N2").Value = 1.076721932794431
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.072131569787574
$ws.Range("D3").Value = 1.072956551518532
$ws.Range("E3").Value = 1.072111763763225
$ws.Range("F3").Value = 1.083832384195884
$ws.Range("I3").Value = 1.049905366561712
$ws.Range("J3").Value = 1.076715272861319
$ws.Range("K3").Value = 1.075468545688602
$ws.Range("L3").Value = 1.074625843368958
$ws.Range("M3").Value = 1.086317842483198
$ws.Range("N3").Value = 1.07824433204975
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.073336157406555
$ws.Range("D4").Value = 1.073915047903025
$ws.Range("E4").Value = 1.073115915354505
$ws.Range("F4").Value = 1.084903340051445
$ws.Range("I4").Value = 1.050242886565995
$ws.Range("J4").Value = 1.077695154419192
$ws.Range("K4").Value = 1.076306765006487
$ws.Range("L4").Value = 1.075509504570381
$ws.Range("M4").Value = 1.087269614830596
$ws.Range("N4").Value = 1.079225605151829
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.073841739359419
$ws.Range("D5").Value = 1.07431722267073
$ws.Range("E5").Value = 1.073537132982751
$ws.Range("F5").Value = 1.08535274855364
$ws.Range("I5").Value = 1.050384039690757
$ws.Range("J5").Value = 1.078106198958757
$ws.Range("K5").Value = 1.076658252646461
$ws.Range("L5").Value = 1.075879949250487
$ws.Range("M5").Value = 1.087668803213616
$ws.Range("N5").Value = 1.079637233421787
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.07392658086902
$ws.Range("D6").Value = 1.074384704428985
$ws.Range("E6").Value = 1.073607803354584
$ws.Range("F6").Value = 1.085428158536999
$ws.Range("I6").Value = 1.050407696693388
$ws.Range("J6").Value = 1.078175162961796
$ws.Range("K6").Value = 1.076717216597246
$ws.Range("L6").Value = 1.075942087627271
$ws.Range("M6").Value = 1.087735774223558
$ws.Range("N6").Value = 1.079706295361622
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.07334291623964
$ws.Range("D7").Value = 1.073920424815839
$ws.Range("E7").Value = 1.073121547307699
$ws.Range("F7").Value = 1.084909348274405
$ws.Range("I7").Value = 1.050244775561111
$ws.Range("J7").Value = 1.077700650322155
$ws.Range("K7").Value = 1.076311465118464
$ws.Range("L7").Value = 1.075514458561031
$ws.Range("M7").Value = 1.087274952464567
$ws.Range("N7").Value = 1.079231108859605
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.070896226415888
$ws.Range("D8").Value = 1.071973178206299
$ws.Range("E8").Value = 1.07108115887175
$ws.Range("F8").Value = 1.082733790942285
$ws.Range("I8").Value = 1.049557492280041
$ws.Range("J8").Value = 1.075709604461412
$ws.Range("K8").Value = 1.07460781782897
$ws.Range("L8").Value = 1.073718115774226
$ws.Range("M8").Value = 1.085340801142321
$ws.Range("N8").Value = 1.077237235485363
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.066556924537604
$ws.Range("D9").Value = 1.068515909396623
$ws.Range("E9").Value = 1.067454775272868
$ws.Range("F9").Value = 1.078872617288866
$ws.Range("I9").Value = 1.048322239727794
$ws.Range("J9").Value = 1.072171191183766
$ws.Range("K9").Value = 1.071575941119528
$ws.Range("L9").Value = 1.070518045877079
$ws.Range("M9").Value = 1.081901401939853
$ws.Range("N9").Value = 1.073693797255002
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.063644070042692
$ws.Range("D10").Value = 1.066192686752266
$ws.Range("E10").Value = 1.065015319177091
$ws.Range("F10").Value = 1.076278909692075
$ws.Range("I10").Value = 1.047482082700026
$ws.Range("J10").Value = 1.069791118424632
$ws.Range("K10").Value = 1.069533768855633
$ws.Range("L10").Value = 1.068360354895853
$ws.Range("M10").Value = 1.07958651093625
$ws.Range("N10").Value = 1.071310344519551
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.062377759455385
$ws.Range("D11").Value = 1.065182141650657
$ws.Range("E11").Value = 1.063953585570877
$ws.Range("F11").Value = 1.075150921322787
$ws.Range("I11").Value = 1.047114244191853
$ws.Range("J11").Value = 1.068755283783097
$ws.Range("K11").Value = 1.068644329746528
$ws.Range("L11").Value = 1.067420062678542
$ws.Range("M11").Value = 1.078578706115828
$ws.Range("N11").Value = 1.070273038874013
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.061906616759626
$ws.Range("D12").Value = 1.064806075133513
$ws.Range("E12").Value = 1.063558373769193
$ws.Range("F12").Value = 1.0747311800889
$ws.Range("I12").Value = 1.046976996959126
$ws.Range("J12").Value = 1.06836972079072
$ws.Range("K12").Value = 1.068313159640631
$ws.Range("L12").Value = 1.067069875943227
$ws.Range("M12").Value = 1.078203525609089
$ws.Range("N12").Value = 1.069886928337963
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.062007714080202
$ws.Range("D13").Value = 1.064886774872742
$ws.Range("E13").Value = 1.06364318623369
$ws.Range("F13").Value = 1.074821250558152
$ws.Range("I13").Value = 1.047006464943476
$ws.Range("J13").Value = 1.068452462236073
$ws.Range("K13").Value = 1.068384232901828
$ws.Range("L13").Value = 1.067145034145443
$ws.Range("M13").Value = 1.078284041287488
$ws.Range("N13").Value = 1.069969787285659
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.062338830672978
$ws.Range("D14").Value = 1.06515107036227
$ws.Range("E14").Value = 1.063920934426843
$ws.Range("F14").Value = 1.075116240926298
$ws.Range("I14").Value = 1.04710291189773
$ws.Range("J14").Value = 1.068723429620351
$ws.Range("K14").Value = 1.068616971393697
$ws.Range("L14").Value = 1.067391135020572
$ws.Range("M14").Value = 1.078547710764022
$ws.Range("N14").Value = 1.070241139474703
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.062542738657062
$ws.Range("D15").Value = 1.065313817569957
$ws.Range("E15").Value = 1.064091952773823
$ws.Range("F15").Value = 1.075297893413162
$ws.Range("I15").Value = 1.047162254233178
$ws.Range("J15").Value = 1.068890273927719
$ws.Range("K15").Value = 1.068760263696353
$ws.Range("L15").Value = 1.067542643358461
$ws.Range("M15").Value = 1.078710054746536
$ws.Range("N15").Value = 1.070408220720124
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.063728003023515
$ws.Range("D16").Value = 1.066259655409875
$ws.Range("E16").Value = 1.065085666735491
$ws.Range("F16").Value = 1.076353665664298
$ws.Range("I16").Value = 1.047506409033854
$ws.Range("J16").Value = 1.069859751140892
$ws.Range("K16").Value = 1.069592687747196
$ws.Range("L16").Value = 1.068422630990377
$ws.Range("M16").Value = 1.079653279324749
$ws.Range("N16").Value = 1.071379074702141
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.064470126079297
$ws.Range("D17").Value = 1.066851717301764
$ws.Range("E17").Value = 1.065707528343475
$ws.Range("F17").Value = 1.077014598496107
$ws.Range("I17").Value = 1.047721199736681
$ws.Range("J17").Value = 1.070466459998974
$ws.Range("K17").Value = 1.070113450723561
$ws.Range("L17").Value = 1.068973004698657
$ws.Range("M17").Value = 1.080243468683089
$ws.Range("N17").Value = 1.071986645156392
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.064902510633174
$ws.Range("D18").Value = 1.067196616842611
$ws.Range("E18").Value = 1.066069726882538
$ws.Range("F18").Value = 1.077399638458381
$ws.Range("I18").Value = 1.04784609377944
$ws.Range("J18").Value = 1.070819838041619
$ws.Range("K18").Value = 1.070416705753548
$ws.Range("L18").Value = 1.069293450983374
$ws.Range("M18").Value = 1.08058719242205
$ws.Range("N18").Value = 1.07234052503639
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.065049861329647
$ws.Range("D19").Value = 1.06731414455273
$ws.Range("E19").Value = 1.066193139159675
$ws.Range("F19").Value = 1.077530847868212
$ws.Range("I19").Value = 1.047888613551735
$ws.Range("J19").Value = 1.070940245757167
$ws.Range("K19").Value = 1.070520024053409
$ws.Range("L19").Value = 1.069402617529032
$ws.Range("M19").Value = 1.080704305083148
$ws.Range("N19").Value = 1.072461103744705
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.064390553455334
$ws.Range("D20").Value = 1.066788240318607
$ws.Range("E20").Value = 1.06564086270266
$ws.Range("F20").Value = 1.076943735521926
$ws.Range("I20").Value = 1.04769819509678
$ws.Range("J20").Value = 1.070401418252392
$ws.Range("K20").Value = 1.07005762933404
$ws.Range("L20").Value = 1.068914014669069
$ws.Range("M20").Value = 1.080180201261719
$ws.Range("N20").Value = 1.071921511043069
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.06224134677208
$ws.Range("D21").Value = 1.065073261535774
$ws.Range("E21").Value = 1.06383916776799
$ws.Range("F21").Value = 1.075029394654819
$ws.Range("I21").Value = 1.047074527732226
$ws.Range("J21").Value = 1.068643658939014
$ws.Range("K21").Value = 1.068548457730646
$ws.Range("L21").Value = 1.067318690008333
$ws.Range("M21").Value = 1.078470089934926
$ws.Range("N21").Value = 1.070161255509848
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.060885538820013
$ws.Range("D22").Value = 1.063990897974751
$ws.Range("E22").Value = 1.062701519657724
$ws.Range("F22").Value = 1.073821387448255
$ws.Range("I22").Value = 1.046678837440549
$ws.Range("J22").Value = 1.067533801856606
$ws.Range("K22").Value = 1.067594986139868
$ws.Range("L22").Value = 1.066310311570356
$ws.Range("M22").Value = 1.077390023111017
$ws.Range("N22").Value = 1.069049822303081
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.061604714316607
$ws.Range("D23").Value = 1.064565072917748
$ws.Range("E23").Value = 1.063305075147704
$ws.Range("F23").Value = 1.074462197524609
$ws.Range("I23").Value = 1.0468889410274
$ws.Range("J23").Value = 1.06812260867195
$ws.Range("K23").Value = 1.068100880748855
$ws.Range("L23").Value = 1.066845384195333
$ws.Range("M23").Value = 1.077963053409858
$ws.Range("N23").Value = 1.06963946529164
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.064426510391643
$ws.Range("D24").Value = 1.06681692419462
$ws.Range("E24").Value = 1.065670987651137
$ws.Range("F24").Value = 1.076975756904526
$ws.Range("I24").Value = 1.047708591107437
$ws.Range("J24").Value = 1.070430809377862
$ws.Range("K24").Value = 1.070082854153122
$ws.Range("L24").Value = 1.068940671507853
$ws.Range("M24").Value = 1.080208790703959
$ws.Range("N24").Value = 1.071950943907309
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.067682170979499
$ws.Range("D25").Value = 1.069412867491277
$ws.Range("E25").Value = 1.06839605867154
$ws.Range("F25").Value = 1.079874197958879
$ws.Range("I25").Value = 1.048644487803738
$ws.Range("J25").Value = 1.073089605919853
$ws.Range("K25").Value = 1.07236337743834
$ws.Range("L25").Value = 1.071349552262397
$ws.Range("M25").Value = 1.082794365580376
$ws.Range("N25").Value = 1.074613516245357
